$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections in column A (attraction names) ---
# Order matters for shared-string append order.
$ws.Range("A10").Value = "TRANSFORMERS: The Ride-3D"
$ws.Range("A7").Value = "MEN IN BLACK Alien Attack"
$ws.Range("A4").Value = "Hogwarts Express - King's Cross Station"

# --- New Courier-New styled column (H9:H15) ---
$cell = $ws.Cells.Item(9, 8)
$cell.Font.Name = "Courier New"
$cell.Font.Size = 14
$cell.Font.Color = 0

$cell.Copy()
$ws.Range("H10:H15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("9:15").RowHeight = 19

# --- Update selection to match the new active area ---
$ws.Range("G5:P26").Select()
